# Apply crypto price/volume updates for Sun Sep 10 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of cell -> new text value (all cells hold plain text, incl. numeric-looking prices)
$updates = [ordered]@{
    'D2' = '26.002.91'
    'E2' = '  -0.13%  '
    'D3' = '1.634.80'
    'E3' = '  -0.57%  '
    'E4' = '  +0.05%  '
    'D5' = '214.08'
    'E5' = '  -1.13%  '
    'E6' = '  -0.64%  '
    'E7' = '  +0.09%  '
    'D8' = '0.251'
    'E8' = '  -2.27%  '
    'E9' = '  -2.45%  '
    'E10' = '  -5.78%  '
    'E11' = '  -0.55%  '
    'D12' = '1.861.78'
    'E12' = '  -0.53%  '
    'E13' = '  -2.05%  '
    'D14' = '1.632.00'
    'D15' = '0.531'
    'E15' = '  -2.62%  '
    'B16' = 'WrappedBTC'
    'C16' = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
    'D16' = '26.011.53'
    'E16' = '  +0.24%  '
    'B17' = 'ShibaInu'
    'C17' = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    'D17' = '0.0₃0746'
    'E17' = '  -2.88%  '
    'D18' = '61.85'
    'E18' = '  -2.01%  '
    'E19' = '  +0.01%  '
    'D20' = '190.61'
    'E20' = '  -1.38%  '
    'E21' = '  -2.31%  '
    'E22' = '  -3.62%  '
    'D23' = '6.15'
    'E23' = '  -1.98%  '
    'E24' = '  +0.90%  '
    'D25' = '143.31'
    'E25' = '  -1.04%  '
    'E27' = '  -2.11%  '
    'D28' = '6.77'
    'E28' = '  -2.61%  '
    'E29' = '  -2.36%  '
    'E30' = '  -1.51%  '
    'D31' = '0.0484'
    'E31' = '  -3.48%  '
    'E32' = '  -2.83%  '
    'E34' = '  -1.64%  '
    'E35' = '  -2.37%  '
    'E36' = '  -3.74%  '
    'D37' = '1.134.03'
    'E37' = '  -0.14%  '
    'D38' = '0.527'
    'E38' = '  -3.18%  '
    'E39' = '  -1.72%  '
    'E40' = '  -1.68%  '
    'D41' = '98.60'
    'E41' = '  -1.15%  '
    'E42' = '  -1.76%  '
    'E43' = '  -4.49%  '
    'D44' = '1.771.97'
    'E44' = '  -0.46%  '
    'E45' = '  -0.79%  '
    'E46' = '  -2.82%  '
    'E47' = '  -0.63%  '
    'D48' = '1.49'
    'E48' = '  +1.40%  '
    'E49' = '  -0.48%  '
    'D50' = '7.52'
    'E50' = '  -2.99%  '
    'D51' = '1.01'
    'E51' = '  -0.02%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings (e.g. "26.002.91") are not
    # coerced into numbers and lose their original formatting/precision.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
